$wb = $excel.ActiveWorkbook

# --- Yearly sheet: update Suzie's Roth IRA October dividend ---
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("F12").Value = 13.08
$wsYearly.Activate() | Out-Null
$wsYearly.Range("D15").Select() | Out-Null

# --- All Time sheet: link 2016 row to Yearly totals via formulas ---
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Range("F7").Formula = "=Yearly!D15"
$wsAllTime.Range("G7").Formula = "=Yearly!E15"
$wsAllTime.Range("H7").Formula = "=Yearly!F15"

# Re-activate "All Time" (it was the active tab before the edit) and update
# its scroll position / selection to match the author's final view.
$wsAllTime.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$wsAllTime.Range("L15").Select() | Out-Null
